$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.817.71'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.290.02'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.56%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.504'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.58'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").Value = '2.646.74'
$ws.Range("D16").Value = '2.292.29'
$ws.Range("E16").Value = '  -0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.772'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").Value = '42.735.31'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("E19").Value = '  -5.25%  '
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("E24").Value = '  -1.48%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '165.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.16%  '
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.66%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.20%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.78'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.09'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.51%  '
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0685'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.81%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.110'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.73'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("D43").Value = '2.007.51'
$ws.Range("E43").Value = '  +0.56%  '
$ws.Range("E44").Value = '  -2.54%  '
$ws.Range("E45").Value = '  -0.89%  '
$ws.Range("E46").Value = '  +3.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.52%  '
$ws.Range("E48").Value = '  -2.38%  '
$ws.Range("D49").Value = '2.514.10'
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.88'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.19%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.07%  '
